$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated price (D) and volume-change (E) values for the cryptos list.
# D-column values that look like plain numbers (e.g. "24.00") would otherwise be
# auto-converted to numeric cells by Excel, so we force text via NumberFormat "@"
# and then restore the default "Normal" style so no stray formatting is left behind.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.433.71'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.20%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.570.75'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.58%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '211.93'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.35%  '
$ws.Range("E6").Value = '  -0.91%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '46.50'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +5.48%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '24.00'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.93%  '
$ws.Range("E11").Value = '  -1.86%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0888'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.01%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.794.99'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.60%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.566.49'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.76%  '
$ws.Range("E15").Value = '  -2.51%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '28.415.21'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.07%  '
$ws.Range("E17").Value = '  -2.30%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '62.08'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.75%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '227.78'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.32%  '
$ws.Range("E20").Value = '  -2.42%  '
$ws.Range("E21").Value = '  -3.09%  '
$ws.Range("E22").Value = '  +0.01%  '
$ws.Range("E23").Value = '  -5.71%  '
$ws.Range("E24").Value = '  -2.73%  '
$ws.Range("E25").Value = '  +7.83%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '150.51'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.00%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '14.95'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.88%  '
$ws.Range("E28").Value = '  -2.38%  '
$ws.Range("E29").Value = '  -3.75%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0481'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.96%  '
$ws.Range("E32").Value = '  -4.07%  '
$ws.Range("E33").Value = '  -1.56%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.08'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.92%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.392.15'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.45%  '
$ws.Range("E36").Value = '  +0.70%  '
$ws.Range("E37").Value = '  -3.49%  '
$ws.Range("E38").Value = '  +0.47%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.63'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.45%  '
$ws.Range("E40").Value = '  -0.68%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.526'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.73%  '
$ws.Range("E42").Value = '  +0.03%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.91'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.83%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.788'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.37%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.974'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.51%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '5.45'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.41%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '62.76'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.97%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.707.28'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.57%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '86.16'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.32%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0₆0103'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.55%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0516'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.94%  '
